$wb = $excel.ActiveWorkbook

# The "Add Panels" sheet contains the 40V load test data table.
$ws1 = $wb.Worksheets.Item("Add Panels")
$ws1.Activate()

# Row 11 ("P405D") is obsolete test data and was removed; deleting the
# entire row shifts rows 12-15 up to become the new rows 11-14.
$ws1.Rows.Item(11).Delete() | Out-Null

# Leave selection on A13, matching the saved worksheet state.
$ws1.Range("A13").Select() | Out-Null
